$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.167
$ws.Range("D3").Value = -0.167
$ws.Range("G2").Value = 0.008064516129032258
$ws.Range("G3").Value = 0.008064516129032258
$ws.Range("H2").Value = 0.008064516129032258
$ws.Range("H3").Value = 0.008064516129032258
$ws.Range("I2").Value = -0.05241935483870968
$ws.Range("I3").Value = -0.05241935483870968
$ws.Range("J2").Value = -0.05241935483870968
$ws.Range("J3").Value = -0.05241935483870968
$ws.Range("K2").Value = -3.27
$ws.Range("K3").Value = -3.27
$ws.Range("L2").Value = -1.318548387096774
$ws.Range("L3").Value = -1.318548387096774
$ws.Range("U2").Value = 1.23
$ws.Range("U3").Value = 1.23
$ws.Range("V2").Value = 0.06473684210526316
$ws.Range("V3").Value = 0.06473684210526316
$ws.Range("W2").Value = 1.003067484662577
$ws.Range("W3").Value = 1.003067484662577
$ws.Range("X2").Value = 0.05995965196265951
$ws.Range("X3").Value = 0.05995965196265951
$ws.Range("Y2").Value = 0.9431078326999172
$ws.Range("Y3").Value = 0.9431078326999172
$ws.Range("Z2").Value = 0.1954910925429607
$ws.Range("Z3").Value = 0.1954910925429607
$ws.Range("AA2").Value = -0.01024751694781649
$ws.Range("AA3").Value = -0.01024751694781649
$ws.Range("AB2").Value = 0.03782406296860871
$ws.Range("AB3").Value = 0.03782406296860871
$ws.Range("AC2").Value = -0.0480715799164252
$ws.Range("AC3").Value = -0.0480715799164252
$ws.Range("AD2").Value = 25.9
$ws.Range("AD3").Value = 25.9
$ws.Range("AE2").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF2").Value = 25.9
$ws.Range("AF3").Value = 25.9
$ws.Range("AG2").Value = 24.67
$ws.Range("AG3").Value = 24.67
$ws.Range("AH2").Value = 0.576837416481069
$ws.Range("AH3").Value = 0.576837416481069
$ws.Range("AI2").Value = 1.001934235976789
$ws.Range("AI3").Value = 1.001934235976789
$ws.Range("AJ2").Value = 0.5649187084955346
$ws.Range("AJ3").Value = 0.5649187084955346
$ws.Range("AK2").Value = 1.002030869212023
$ws.Range("AK3").Value = 1.002030869212023
$ws.Range("AL2").Value = 0.917
$ws.Range("AL3").Value = 0.917
$ws.Range("AM2").Value = 0.2190000000000001
$ws.Range("AM3").Value = 0.2190000000000001
$ws.Range("AN2").Value = -25900
$ws.Range("AN3").Value = -25900
$ws.Range("AO2").Value = -0.1417666303162486
$ws.Range("AO3").Value = -0.1417666303162486
$ws.Range("AP2").Value = -24670
$ws.Range("AP3").Value = -24670
$ws.Range("AQ2").Value = -0.5936073059360728
$ws.Range("AQ3").Value = -0.5936073059360728
